$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

function Set-CellText([int]$row, [int]$col, [string]$value) {
    $t.Cell($row, $col).Range.Text = $value
}

Set-CellText 1 1 "13÷7=1, 6"
Set-CellText 1 2 "36÷9=4, 0"
Set-CellText 1 3 "34÷5=6, 4"
Set-CellText 1 4 "21÷2=10, 1"
Set-CellText 1 5 "36÷4=9, 0"
Set-CellText 5 1 "42÷9=4, 6"
Set-CellText 5 2 "65÷5=13, 0"
Set-CellText 5 3 "58÷4=14, 2"
Set-CellText 5 4 "62÷6=10, 2"
Set-CellText 5 5 "17÷4=4, 1"
Set-CellText 9 1 "46÷9=5, 1"
Set-CellText 9 2 "42÷6=7, 0"
Set-CellText 9 3 "72÷4=18, 0"
Set-CellText 9 4 "72÷8=9, 0"
Set-CellText 9 5 "62÷9=6, 8"
Set-CellText 13 1 "85÷5=17, 0"
Set-CellText 13 2 "17÷8=2, 1"
Set-CellText 13 3 "98÷2=49, 0"
Set-CellText 13 4 "81÷9=9, 0"
Set-CellText 13 5 "24÷3=8, 0"
Set-CellText 17 1 "58÷6=9, 4"
Set-CellText 17 2 "97÷9=10, 7"
Set-CellText 17 3 "92÷5=18, 2"
Set-CellText 17 4 "14÷3=4, 2"
Set-CellText 17 5 "78÷8=9, 6"

Write-Output "updated $($t.Rows.Count) x $($t.Columns.Count) table"
